# Update FedEx shipment tracking numbers (column P, "ShipmentTracking") for
# rows 2-25 with the new tracking numbers generated on 1st April 2022.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$trackingNumbers = @{
    2  = "320018208097"
    3  = "320018208101"
    4  = "320018208134"
    5  = "320018208156"
    6  = "320018208190"
    7  = "320018208215"
    8  = "320018208248"
    9  = "320018208260"
    10 = "320018208292"
    11 = "320018208318"
    12 = "320018208351"
    13 = "320018208373"
    14 = "320018208400"
    15 = "320018208421"
    16 = "320018208454"
    17 = "320018208476"
    18 = "320018208513"
    19 = "320018208535"
    20 = "320018208568"
    21 = "320018208580"
    22 = "320018208616"
    23 = "320018208627"
    24 = "320018208638"
    25 = "320018208649"
}

foreach ($row in $trackingNumbers.Keys) {
    $ws.Cells.Item($row, 16).Value = $trackingNumbers[$row]
}
